# taskboard.xlsx edit script
# Adds a "Descricao" text for each existing task row, inserts a new task
# row ("Unificar as bases de dados de Wine"), then applies the
# alignment/format/width changes and resizes the table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Fill in the "Descricao" column (C) for the existing rows, top to
#    bottom, BEFORE inserting the new row - this reproduces the shared
#    string insertion order used by the original author.
$ws.Range("C6").Value  = "Verificar se as tabelas tem dado faltante, quantificar por variavel (%). Avaliar a faixa dinamica de cada variavel."
$ws.Range("C7").Value  = "Histogramas com as distribuicoes das variaveis. Utilizar tambem KDE. Analisar padrões existentes"
$ws.Range("C8").Value  = "Histogramas e KDE para diferentes grupos. Classificao sao as classes. Regressao, binarizar a variavel de saida para formar os grupos. Eventualmente variaveis que nao sejam continuas"
$ws.Range("C9").Value  = "Calcular os coeficientes de correlacao e scatter-plot."
$ws.Range("C10").Value = "Calcular os box-plot para cada variavel no mesmo graficao (se aplicavel): identificar outliers"
$ws.Range("C11").Value = "Plot-scatter para cada variavel (pandas)"
$ws.Range("C12").Value = "Reefetuar as analises com a normalizacao"
$ws.Range("C13").Value = "Reefetuar as analises com a normalizacao"
$ws.Range("C14").Value = "Reefetuar as analises com a normalizacao (onde aplicavel)"

# 2) Insert a new row at row 6 (pushes the rows filled in above down by
#    one, to rows 7-15) and fill it with the new "Wine" unification task.
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = "Data Praparation"
$ws.Range("B6").Value = "Unificar as bases de dados de Wine"
$ws.Range("C6").Value = "Unificar as bases de vinho branco e tinto e adicionar a informacao branco/tinto como variavel (coluna)"
$ws.Range("D6").Value = "Dados preparados no ambiente"
$ws.Range("E6").Value = "Notebook fazendo a avaliação"

# 3) Grow the table to include the new row.
$lo.Resize($ws.Range("A1:F15"))

# 4) Formatting: header C1 gets wrap text, the whole data body gets
#    vertical-center alignment, and column C (Descricao) gets wrap text.
$ws.Range("C1").WrapText = $true
$ws.Range("A2:F15").VerticalAlignment = -4108
$ws.Range("C2:C15").WrapText = $true

# 5) Column C width (matches the widened "Descricao" column).
$ws.Columns("C").ColumnWidth = 42.57

# 6) Row heights for the wrapped description rows (auto-fit result).
$ws.Rows(6).RowHeight = 43.5
$ws.Rows(7).RowHeight = 43.5
$ws.Rows(8).RowHeight = 29
$ws.Rows(9).RowHeight = 58
$ws.Rows(10).RowHeight = 29
$ws.Rows(11).RowHeight = 29
$ws.Rows(15).RowHeight = 29

# 7) Selection moves to A6, matching the saved view state.
[void]$ws.Range("A6").Select()
